$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Cells.Item(3, 2).Value = 0.01
$ws.Cells.Item(6, 2).Value = 44516.99100015465
$ws.Cells.Item(7, 2).Value = 11110450.29222555
$ws.Cells.Item(8, 2).Value = 21894283.74316604
$ws.Cells.Item(10, 2).Value = 4015769.133974414

# --- Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Cells.Item(2, 2).Value = 81828.8692115825
$ws.Cells.Item(2, 3).Value = 86920.14920292264
$ws.Cells.Item(2, 4).Value = 87500.02638886587
$ws.Cells.Item(2, 6).Value = 89830.69576735962
$ws.Cells.Item(2, 7).Value = 89830.6957673596
$ws.Cells.Item(2, 8).Value = 89830.69576735958
$ws.Cells.Item(2, 9).Value = 89830.69576735958
$ws.Cells.Item(2, 10).Value = 89830.69576735959
$ws.Cells.Item(2, 12).Value = 89830.69576735958
$ws.Cells.Item(2, 13).Value = 89830.69576735963
$ws.Cells.Item(2, 14).Value = 89830.69576735962
$ws.Cells.Item(2, 15).Value = 89830.69576735958
$ws.Cells.Item(2, 16).Value = 89830.69576735959
$ws.Cells.Item(3, 2).Value = 22175.99473096512
$ws.Cells.Item(3, 3).Value = 172096.213521813
$ws.Cells.Item(3, 4).Value = 19504.79678339268
$ws.Cells.Item(4, 2).Value = 66348.57580579144
$ws.Cells.Item(4, 3).Value = 33232.85497694127
$ws.Cells.Item(5, 2).Value = 34129.29174257201
$ws.Cells.Item(5, 3).Value = 38249.31064740147
$ws.Cells.Item(6, 2).Value = -55895.77936691431
$ws.Cells.Item(6, 3).Value = -171389.597576312
$ws.Cells.Item(6, 4).Value = -14656.26526343197
$ws.Cells.Item(6, 5).Value = -20264.66179571922
$ws.Cells.Item(6, 6).Value = 53277.39439424062
$ws.Cells.Item(6, 7).Value = 53277.3943942406
$ws.Cells.Item(6, 8).Value = 53277.39439424057
$ws.Cells.Item(6, 9).Value = 53277.39439424057
$ws.Cells.Item(6, 10).Value = 53277.39439424059
$ws.Cells.Item(6, 11).Value = 53277.39439424057
$ws.Cells.Item(6, 12).Value = 53277.39439424057
$ws.Cells.Item(6, 13).Value = 53277.39439424063
$ws.Cells.Item(6, 14).Value = 53277.39439424062
$ws.Cells.Item(6, 15).Value = 53277.39439424057
$ws.Cells.Item(6, 16).Value = 53277.39439424059

# --- Installed Capacities ---
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Cells.Item(3, 2).Value = 23.01338268678932
$ws.Cells.Item(3, 3).Value = 212.005075568875

# --- Added Capacities ---
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Cells.Item(3, 2).Value = 23.01338268678932
$ws.Cells.Item(3, 3).Value = 188.9916928820857
$ws.Cells.Item(3, 4).Value = 22.7749625571192

# --- PV Dispatch ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Cells.Item(2, 7).Value = 0.0925161113036756
$ws.Cells.Item(2, 8).Value = 0.9474806248887678
$ws.Cells.Item(2, 9).Value = 3.566727381034957
$ws.Cells.Item(2, 10).Value = 7.852189301760343
$ws.Cells.Item(2, 11).Value = 11.76839629324493
$ws.Cells.Item(2, 12).Value = 14.5997362345548
$ws.Cells.Item(2, 13).Value = 16.24501962895154
$ws.Cells.Item(2, 14).Value = 16.50788103019311
$ws.Cells.Item(2, 15).Value = 15.58792394841718
$ws.Cells.Item(2, 16).Value = 13.30393245060769
$ws.Cells.Item(2, 17).Value = 9.990699214544804
$ws.Cells.Item(2, 18).Value = 5.811515176679517
$ws.Cells.Item(2, 19).Value = 2.10821088633251
$ws.Cells.Item(2, 20).Value = 0.4049892772318401
$ws.Cells.Item(2, 21).Value = 0.007401288904294046
$ws.Cells.Item(3, 7).Value = 0.04950048351498081
$ws.Cells.Item(3, 8).Value = 0.4780704592104726
$ws.Cells.Item(3, 9).Value = 1.704292963125436
$ws.Cells.Item(3, 10).Value = 4.676710155246498
$ws.Cells.Item(3, 11).Value = 7.993242550750212
$ws.Cells.Item(3, 12).Value = 10.74790103688213
$ws.Cells.Item(3, 13).Value = 12.54229356430018
$ws.Cells.Item(3, 14).Value = 12.87425075418792
$ws.Cells.Item(3, 15).Value = 11.77742425104019
$ws.Cells.Item(3, 16).Value = 9.452421277522959
$ws.Cells.Item(3, 17).Value = 6.318693299210532
$ws.Cells.Item(3, 18).Value = 3.073372125605564
$ws.Cells.Item(3, 19).Value = 0.9194497705523841
$ws.Cells.Item(3, 20).Value = 0.1995216857467866
$ws.Cells.Item(3, 21).Value = 0.003256610757564528
$ws.Cells.Item(4, 7).Value = 0.04149954254994795
$ws.Cells.Item(4, 8).Value = 0.3689686601259011
$ws.Cells.Item(4, 9).Value = 1.248004425047526
$ws.Cells.Item(4, 10).Value = 2.93401765828132
$ws.Cells.Item(4, 11).Value = 4.821492307166679
$ws.Cells.Item(4, 12).Value = 6.169850171471354
$ws.Cells.Item(4, 13).Value = 6.505241928988658
$ws.Cells.Item(4, 14).Value = 6.350561815847948
$ws.Cells.Item(4, 15).Value = 5.865771705150827
$ws.Cells.Item(4, 16).Value = 5.019181037131885
$ws.Cells.Item(4, 17).Value = 3.475020785705187
$ws.Cells.Item(4, 18).Value = 1.865970340473114
$ws.Cells.Item(4, 19).Value = 0.7232238460750018
$ws.Cells.Item(4, 20).Value = 0.1773162272588685
$ws.Cells.Item(4, 21).Value = 0.002263611411815345
$ws.Cells.Item(5, 7).Value = 0.8522817108296475
$ws.Cells.Item(5, 8).Value = 8.728430071034131
$ws.Cells.Item(5, 9).Value = 32.85759065676002
$ws.Cells.Item(5, 10).Value = 72.33634485452787
$ws.Cells.Item(5, 11).Value = 108.4134296739469
$ws.Cells.Item(5, 12).Value = 134.4964460817497
$ws.Cells.Item(5, 13).Value = 149.6532109567165
$ws.Cells.Item(5, 14).Value = 152.0747563676112
$ws.Cells.Item(5, 15).Value = 143.5998801055489
$ws.Cells.Item(5, 16).Value = 122.5591753694419
$ws.Cells.Item(5, 17).Value = 92.03683660035516
$ws.Cells.Item(5, 18).Value = 53.5371410179029
$ws.Cells.Item(5, 19).Value = 19.42136948553061
$ws.Cells.Item(5, 20).Value = 3.730863189156784
$ws.Cells.Item(5, 21).Value = 0.0681825368663718
$ws.Cells.Item(6, 7).Value = 0.4560109172613537
$ws.Cells.Item(6, 8).Value = 4.40410543776097
$ws.Cells.Item(6, 9).Value = 15.7003758793931
$ws.Cells.Item(6, 10).Value = 43.08303144178431
$ws.Cells.Item(6, 11).Value = 73.63576289829501
$ws.Cells.Item(6, 12).Value = 99.0123704343128
$ws.Cells.Item(6, 13).Value = 115.5427661850369
$ws.Cells.Item(6, 14).Value = 118.6008393977238
$ws.Cells.Item(6, 15).Value = 108.4965974941959
$ws.Cells.Item(6, 16).Value = 87.07808471809132
$ws.Cells.Item(6, 17).Value = 58.20939357883527
$ws.Cells.Item(6, 18).Value = 28.31267782785845
$ws.Cells.Item(6, 19).Value = 8.470202783341367
$ws.Cells.Item(6, 20).Value = 1.838044004224491
$ws.Cells.Item(6, 21).Value = 0.03000071824087855
$ws.Cells.Item(7, 7).Value = 0.3823042346323975
$ws.Cells.Item(7, 8).Value = 3.399032195186227
$ws.Cells.Item(7, 9).Value = 11.49693098330883
$ws.Cells.Item(7, 10).Value = 27.0289093885105
$ws.Cells.Item(7, 11).Value = 44.41680107820035
$ws.Cells.Item(7, 12).Value = 56.83821321071117
$ws.Cells.Item(7, 13).Value = 59.92792652514935
$ws.Cells.Item(7, 14).Value = 58.50297437788319
$ws.Cells.Item(7, 15).Value = 54.03696581876834
$ws.Cells.Item(7, 16).Value = 46.2379594322674
$ws.Cells.Item(7, 17).Value = 32.01276641090012
$ws.Cells.Item(7, 18).Value = 17.18978858628943
$ws.Cells.Item(7, 19).Value = 6.662520161730052
$ws.Cells.Item(7, 20).Value = 1.633481729792971
$ws.Cells.Item(7, 21).Value = 0.02085295825267625
$ws.Cells.Item(11, 9).Value = 50.57440593705229
$ws.Cells.Item(11, 18).Value = 82.40437136225565
$ws.Cells.Item(12, 7).Value = 0.7018920371314099
$ws.Cells.Item(12, 18).Value = 43.57887577908878
$ws.Cells.Item(13, 8).Value = 5.231790602914801
$ws.Cells.Item(13, 12).Value = 87.48538192236053
$ws.Cells.Item(13, 13).Value = 92.24106888145177

# --- Fed-in Capacity ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Cells.Item(2, 11).Value = 208.3214547517356
$ws.Cells.Item(2, 12).Value = 221.1666787354324
$ws.Cells.Item(2, 13).Value = 214.1012135983212
$ws.Cells.Item(2, 14).Value = 212.9051825663978
$ws.Cells.Item(2, 15).Value = 214.5102874732696
$ws.Cells.Item(2, 16).Value = 217.9290633046618
$ws.Cells.Item(3, 10).Value = 122.1609165114202
$ws.Cells.Item(3, 11).Value = 129.8481964236088
$ws.Cells.Item(3, 12).Value = 127.8064787429921
$ws.Cells.Item(3, 13).Value = 129.5917403577182
$ws.Cells.Item(3, 14).Value = 118.4674613291454
$ws.Cells.Item(3, 15).Value = 130.8188201934042
$ws.Cells.Item(3, 16).Value = 124.5219861368073
$ws.Cells.Item(3, 17).Value = 133.663080786811
$ws.Cells.Item(4, 12).Value = 128.7148261097669
$ws.Cells.Item(4, 13).Value = 132.4205420186164
$ws.Cells.Item(4, 14).Value = 121.3349826493852
$ws.Cells.Item(4, 15).Value = 132.590766746692
$ws.Cells.Item(4, 16).Value = 132.7088230120157
$ws.Cells.Item(5, 10).Value = 108.7095596721584
$ws.Cells.Item(5, 11).Value = 111.6764213710337
$ws.Cells.Item(5, 12).Value = 101.2699688882375
$ws.Cells.Item(5, 13).Value = 80.69302227055627
$ws.Cells.Item(5, 14).Value = 77.33830722897972
$ws.Cells.Item(5, 15).Value = 86.49833131613784
$ws.Cells.Item(5, 16).Value = 108.6738203858276
$ws.Cells.Item(5, 17).Value = 130.2688532740943
$ws.Cells.Item(6, 10).Value = 83.75459522488239
$ws.Cells.Item(6, 11).Value = 64.20567607606398
$ws.Cells.Item(6, 12).Value = 39.54200934556138
$ws.Cells.Item(6, 13).Value = 26.59126773698146
$ws.Cells.Item(6, 14).Value = 12.74087268560955
$ws.Cells.Item(6, 15).Value = 34.09964695024857
$ws.Cells.Item(6, 16).Value = 46.89632269623893
$ws.Cells.Item(6, 17).Value = 81.77238050718626
$ws.Cells.Item(7, 11).Value = 84.59645709233503
$ws.Cells.Item(7, 12).Value = 78.04646307052711
$ws.Cells.Item(7, 13).Value = 78.9978574224557
$ws.Cells.Item(7, 14).Value = 69.18257008735
$ws.Cells.Item(7, 15).Value = 84.41957263307444
$ws.Cells.Item(7, 16).Value = 91.49004461688021
$ws.Cells.Item(11, 15).Value = 9.069265482343809
$ws.Cells.Item(11, 17).Value = 80.64258426171739
$ws.Cells.Item(13, 12).Value = 47.39929435887775
$ws.Cells.Item(13, 13).Value = 46.68471506615329

# --- Unmet Demand ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Cells.Item(2, 7).Value = 415.2102214038314
$ws.Cells.Item(2, 8).Value = 338.5273214908784
$ws.Cells.Item(2, 9).Value = 206.909162189371
$ws.Cells.Item(2, 10).Value = 4.097100052852174
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 144.0576027644702
$ws.Cells.Item(2, 19).Value = 206.9118586999128
$ws.Cells.Item(2, 20).Value = 222.6908602868995
$ws.Cells.Item(2, 21).Value = 251.3382516189322
$ws.Cells.Item(3, 7).Value = 137.2940166796957
$ws.Cells.Item(3, 8).Value = 111.757373777286
$ws.Cells.Item(3, 9).Value = 87.69233988828964
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 18).Value = 97.08446202703757
$ws.Cells.Item(3, 19).Value = 170.7637213332854
$ws.Cells.Item(3, 20).Value = 199.9652070090748
$ws.Cells.Item(3, 21).Value = 225.9381254702172
$ws.Cells.Item(4, 7).Value = 167.9494798159088
$ws.Cells.Item(4, 8).Value = 161.8582038473137
$ws.Cells.Item(4, 9).Value = 154.2024705022108
$ws.Cells.Item(4, 10).Value = 90.42516245839145
$ws.Cells.Item(4, 11).Value = 17.44799951871617
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 82.6870224659892
$ws.Cells.Item(4, 18).Value = 175.4274210366964
$ws.Cells.Item(4, 19).Value = 223.2933741908973
$ws.Cells.Item(4, 20).Value = 227.7682732010226
$ws.Cells.Item(4, 21).Value = 286.3167657450791
$ws.Cells.Item(5, 7).Value = 414.4504558043054
$ws.Cells.Item(5, 8).Value = 330.746372044733
$ws.Cells.Item(5, 9).Value = 177.6182989136459
$ws.Cells.Item(5, 18).Value = 96.33197692324678
$ws.Cells.Item(5, 19).Value = 189.5987001007147
$ws.Cells.Item(5, 20).Value = 219.3649863749746
$ws.Cells.Item(5, 21).Value = 251.2774703709701
$ws.Cells.Item(6, 7).Value = 136.8875062459493
$ws.Cells.Item(6, 8).Value = 107.8313387987355
$ws.Cells.Item(6, 9).Value = 73.69625697202197
$ws.Cells.Item(6, 18).Value = 71.84515632478468
$ws.Cells.Item(6, 19).Value = 163.2129683204965
$ws.Cells.Item(6, 20).Value = 198.3266846905971
$ws.Cells.Item(6, 21).Value = 225.9113813627339
$ws.Cells.Item(7, 7).Value = 167.6086751238264
$ws.Cells.Item(7, 8).Value = 158.8281403122533
$ws.Cells.Item(7, 9).Value = 143.9535439439495
$ws.Cells.Item(7, 10).Value = 66.33027072816228
$ws.Cells.Item(7, 17).Value = 54.14927684079426
$ws.Cells.Item(7, 18).Value = 160.10360279088
$ws.Cells.Item(7, 19).Value = 217.3540778752422
$ws.Cells.Item(7, 20).Value = 226.3121076984885
$ws.Cells.Item(7, 21).Value = 286.2981763982382
$ws.Cells.Item(11, 18).Value = 67.46474657889404
$ws.Cells.Item(12, 18).Value = 56.57895837355435

# --- Household Surplus ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Cells.Item(2, 2).Value = 366387.7629149273
$ws.Cells.Item(3, 2).Value = 351940.0687209912
